# Update the "Förändrad" (Changed) date column for all existing data rows
# (row 2 through row 381) from 2023-09-23 (45192) to 2023-10-03 (45202).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C381").Value2 = 45202

# Row 381 gets an explicit row height (matches the new trailing rows).
$ws.Rows.Item(381).RowHeight = 15

# Append a new data row (row 382) for case "A 45766-2023".
$ws.Range("A382").Value2 = "A 45766-2023"

$ws.Range("B382").Value2 = 45195
$ws.Range("B382").NumberFormat = "YYYY-MM-DD"

$ws.Range("C382").Value2 = 45202
$ws.Range("C382").NumberFormat = "YYYY-MM-DD"

$ws.Range("D382").Value2 = "SÖDERMANLANDS LÄN"
$ws.Range("E382").Value2 = "GNESTA"

$ws.Range("G382").Value2 = 6.7
$ws.Range("H382").Value2 = 0
$ws.Range("I382").Value2 = 0
$ws.Range("J382").Value2 = 0
$ws.Range("K382").Value2 = 0
$ws.Range("L382").Value2 = 0
$ws.Range("M382").Value2 = 0
$ws.Range("N382").Value2 = 0
$ws.Range("O382").Value2 = 0
$ws.Range("P382").Value2 = 0
$ws.Range("Q382").Value2 = 0

# R382 stays an empty (wrapped) cell, matching the style used throughout column R.
$ws.Range("R382").WrapText = $true
